$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 246: correct a tiny floating point rounding difference in the
# stored timestamp (no other change to that row). ---
$ws.Range("B246").Value = 44503.68519162037

# --- Rows 247-281: append 35 new WPM test-run entries. ---
# Columns: A=Nr. B=Datum(serial) C=WPM D=Genauigkeit(text) E=Tastenanschlaege
# (gesamt) F=(richtig) G=(falsch) H=Korrekte Woerter I=Falsche Woerter
$newRows = @(
    @(247, 246, 44511.51274792824, 78, "96.05", 398, 389, 9, 53, 1),
    @(248, 247, 44515.56589346065, 82, "96.48", 421, 411, 10, 57, 1),
    @(249, 248, 44517.71932373843, 77, "97.46", 391, 384, 7, 53, 1),
    @(250, 249, 44519.4644558912, 67, "91.08", 357, 337, 20, 47, 2),
    @(251, 250, 44519.46520854167, 76, "95.26", 382, 382, 0, 53, 0),
    @(252, 251, 44519.46599408565, 69, "96.11", 346, 346, 0, 46, 0),
    @(253, 252, 44522.53524550926, 68, "87.89", 375, 341, 34, 46, 4),
    @(254, 253, 44522.53630471065, 71, "95.43", 362, 355, 7, 53, 1),
    @(255, 254, 44522.53792792824, 77, "95.56", 401, 387, 14, 52, 2),
    @(256, 255, 44523.6913862037, 76, "98.2", 382, 382, 0, 54, 0),
    @(257, 256, 44523.6927865625, 76, "90.89", 405, 379, 26, 52, 2),
    @(258, 257, 44523.69372001157, 71, "93.14", 368, 353, 15, 48, 2),
    @(259, 258, 44523.69758021991, 74, "95.1", 376, 369, 7, 56, 1),
    @(260, 259, 44526.40262489583, 74, "96.84", 368, 368, 0, 51, 0),
    @(261, 260, 44526.40426511574, 83, "96.05", 426, 413, 13, 54, 2),
    @(262, 261, 44526.40581282407, 73, "91.5", 395, 366, 29, 49, 4),
    @(263, 262, 44529.68288987268, 68, "88.6", 371, 342, 29, 45, 4),
    @(264, 263, 44529.68367783564, 79, "94.5", 411, 395, 16, 57, 2),
    @(265, 264, 44529.68447078704, 81, "96.9", 414, 407, 7, 57, 1),
    @(266, 265, 44529.68553372685, 79, "96.57", 398, 394, 4, 53, 1),
    @(267, 266, 44531.51940063657, 76, "93.61", 396, 381, 15, 48, 2),
    @(268, 267, 44538.53557773148, 74, "91.38", 393, 371, 22, 49, 3),
    @(269, 268, 44538.53750008102, 76, "93.86", 388, 382, 6, 52, 1),
    @(270, 269, 44538.53843518518, 84, "94.16", 443, 419, 24, 60, 3),
    @(271, 270, 44539.41822783565, 75, "93.78", 391, 377, 14, 54, 2),
    @(272, 271, 44539.56732179398, 78, "92.16", 411, 388, 23, 55, 3),
    @(273, 272, 44543.5055310301, 70, "93.1", 358, 351, 7, 52, 1),
    @(274, 273, 44543.51116302083, 82, "96.04", 424, 412, 12, 57, 1),
    @(275, 274, 44543.51307186342, 68, "86.51", 379, 340, 39, 47, 3),
    @(276, 275, 44543.51389138889, 77, "94.15", 402, 386, 16, 55, 2),
    @(277, 276, 44572.43762675926, 73, "94.3", 372, 364, 8, 49, 1),
    @(278, 277, 44572.46692445602, 72, "90.95", 385, 362, 23, 51, 3),
    @(279, 278, 44572.46776936343, 82, "98.56", 410, 410, 0, 59, 0),
    @(280, 279, 44572.47360372685, 80, "92.99", 416, 398, 18, 55, 2),
    @(281, 280, 44578.64101636007, 79, "93.19", 419, 397, 22, 53, 2)
)

# Step 1: write column D (Genauigkeit) as literal text first, forcing a
# "@" text format so the numeric-looking strings ("96.05" etc.) are not
# auto-converted to numbers, matching the source data.
foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Step 2: copy the style/format block from row 246 (the last populated
# data row) down across all of the newly-populated rows in one shot. This
# overwrites the "@" format set above with the correct cell style while
# leaving the text we just entered as text.
$ws.Range("A246:I246").Copy()
$ws.Range("A247:I281").PasteSpecial(-4122)

# Step 3: write the remaining numeric / date columns for every new row.
foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}
